{"js": "// Fix register order in instructions.\n// Three \"op ...\" operand-format example lines list operands in the wrong\n// order; reorder them so that `destination` precedes the source operand(s),\n// matching the actual instruction-format tables above each line.\n\nconst body = context.document.body;\n\nconst replacements = [\n  {\n    find: \"op source1, source2, destination\",\n    replace: \"op destination, source1, source2\",\n  },\n  {\n    find: \"op source, destination, shiftAmount\",\n    replace: \"op destination, source, shiftAmount\",\n  },\n  {\n    find: \"op destination, source, offset\",\n    replace: \"op source, destination, offset\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const found = body.search(find, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${find}\"`);\n  }\n\n  found.items[0].insertText(replace, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Fix register order in instructions.\n# Three \"op ...\" operand-format example lines list operands in the wrong\n# order; reorder them so that `destination` precedes the source operand(s),\n# matching the actual instruction-format tables above each line.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"op source1, source2, destination\"; Replace = \"op destination, source1, source2\" },\n    @{ Find = \"op source, destination, shiftAmount\"; Replace = \"op destination, source, shiftAmount\" },\n    @{ Find = \"op destination, source, offset\"; Replace = \"op source, destination, offset\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $found = $range.Find.Execute(\n        $r.Find,    # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $r.Replace, # ReplaceWith\n        2           # Replace (wdReplaceOne)\n    )\n    if (-not $found) {\n        throw \"Could not find text to replace: $($r.Find)\"\n    }\n}\n"}
